# Update "想去人数" (want-to-go count) figures in column F of the
# "展览" and "全部类型" worksheets, reflecting freshly generated site data.

$wb = $excel.ActiveWorkbook

# -- 展览 (Exhibition) sheet --
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 182
$wsExhibit.Range("F5").Value = 302
$wsExhibit.Range("F7").Value = 253
$wsExhibit.Range("F8").Value = 2341
$wsExhibit.Range("F10").Value = 5885
$wsExhibit.Range("F11").Value = 146

# -- 全部类型 (All types) sheet --
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 182
$wsAll.Range("F6").Value = 302
$wsAll.Range("F8").Value = 253
$wsAll.Range("F11").Value = 2341
$wsAll.Range("F13").Value = 5885
$wsAll.Range("F14").Value = 146
